$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace an entire paragraph (identified by a distinctive substring
# of its current text) with a brand-new <w:p> ... </w:p> fragment, using
# Range.InsertXML so we get full control over pPr/rPr instead of fighting
# the object model's per-property quirks.
# ---------------------------------------------------------------------------
function Set-ParagraphXml {
    param(
        [string]$MatchText,
        [string]$InnerParagraphXml
    )

    $target = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($MatchText)) {
            $target = $p
            break
        }
    }
    if ($target -eq $null) {
        throw "Paragraph containing '$MatchText' not found"
    }

    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $InnerParagraphXml + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $target.Range.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# 1) "We strive to..." -> "The self funding mechanism ..."
# ---------------------------------------------------------------------------
$para1 = '<w:p>' +
    '<w:pPr>' +
        '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
            '<w:b/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Open Sans" w:hAnsi="Open Sans" w:cs="Open Sans" w:eastAsia="Times New Roman"/>' +
            '<w:b/>' +
            '<w:color w:val="000000"/>' +
            '<w:sz w:val="21"/>' +
            '<w:szCs w:val="21"/>' +
            '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
        '</w:rPr>' +
        '<w:t xml:space="preserve">The self funding mechanism that pays for development and encourages the community to decide how budget funds are spent.</w:t>' +
    '</w:r>' +
'</w:p>'

Set-ParagraphXml "We strive to" $para1

# ---------------------------------------------------------------------------
# 2) "SmartCash mining prevents ..." -> "SmartCash now has SmartMining ..."
# ---------------------------------------------------------------------------
$para2 = '<w:p>' +
    '<w:pPr>' +
        '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma" w:eastAsia="Times New Roman"/>' +
            '<w:color w:val="252525"/>' +
            '<w:sz w:val="28"/>' +
            '<w:szCs w:val="28"/>' +
            '<w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/>' +
        '</w:rPr>' +
        ('<w:t xml:space="preserve">SmartCash now has SmartMining that prevents mining attacks. Mining can be done by anyone with a computer with one or more graphics cards. ASICs have yet to be created for the Keccak mining algorithm and it' + [char]0x2019 + 's probably safe to assume no ASICs will be created for quite some time.</w:t>') +
    '</w:r>' +
'</w:p>'

Set-ParagraphXml "SmartCash mining prevents" $para2

Write-Output "edits applied"
